$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'60.199.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.24%  '
$ws.Range("D3").Formula = "'2.582.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Formula = "'506.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.38%  '
$ws.Range("D6").Formula = "'153.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.44%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").Value = '  -8.05%  '
$ws.Range("D9").Formula = "'2.585.72"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("D10").Formula = "'6.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.78%  '
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("E12").Value = '  +1.54%  '
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("D14").Formula = "'3.034.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Formula = "'60.163.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.33%  '
$ws.Range("E16").Value = '  -1.45%  '
$ws.Range("D17").Formula = "'0.0000140"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.69%  '
$ws.Range("D18").Formula = "'2.584.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("E19").Value = '  +1.76%  '
$ws.Range("D20").Formula = "'345.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.11%  '
$ws.Range("D21").Formula = "'10.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.76%  '
$ws.Range("D22").Formula = "'6.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.09%  '
$ws.Range("E23").Value = '  -0.90%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Formula = "'0.420"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.57%  '
$ws.Range("D26").Formula = "'0.166"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.98%  '
$ws.Range("D27").Formula = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("D28").Formula = "'0.0₃0842"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.41%  '
$ws.Range("D29").Formula = "'7.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.53%  '
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").Formula = "'19.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").Formula = "'153.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.59%  '
$ws.Range("E33").Value = '  -0.87%  '
$ws.Range("D34").Formula = "'5.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.57%  '
$ws.Range("D35").Formula = "'3.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.50%  '
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("D37").Formula = "'0.858"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +17.71%  '
$ws.Range("E38").Value = '  -0.14%  '
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("E40").Value = '  +1.75%  '
$ws.Range("D41").Formula = "'35.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.23%  '
$ws.Range("D42").Formula = "'294.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.49%  '
$ws.Range("D43").Formula = "'0.0998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.84%  '
$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D44").Formula = "'0.0559"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.87%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Formula = "'0.613"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.74%  '
$ws.Range("D46").Formula = "'0.996"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.40%  '
$ws.Range("D47").Formula = "'19.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.36%  '
$ws.Range("D48").Formula = "'4.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.34%  '
$ws.Range("E49").Value = '  -0.67%  '
$ws.Range("E50").Value = '  -0.28%  '
$ws.Range("D51").Formula = "'2.000.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.76%  '
